# Orders workbook update:
#  - Insert two new rows (order ids 2 and 3) after the current row 2,
#    pushing the existing rows 3-6 down to rows 5-8.
#  - Populate the two new rows with new client data (C00005 / C00006),
#    whose fund-id column (D) is entered as TEXT ("111" / "222") instead
#    of a number - per commit message "Add trailing Zeros to a string".
#  - Re-point/renumber the "order id" column (A) for every data row.
#  - Turn the L column into a running "H + offset" formula for every row
#    (L2 used to be a hard-coded 1000; row 5 - formerly row 3 - used to be
#    a hard-coded 2020.222).
#  - Leave the final selection on L4, matching the author's last action.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank rows above the old row 3 ------------------------
$ws.Rows("3:4").Insert()
$ws.Rows("3:4").RowHeight = 18

# --- 2. Brand-new rows 3 & 4 (order ids 2 & 3 / C00005 & C00006) --------
# Filled column-by-column (both rows' client_code, then both rows' fund_id,
# then both rows' bank.code, ...) to mirror how the data was pasted in.
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

$ws.Range("B3").Value = "{epoch_id}"
$ws.Range("B4").Value = "{epoch_id}"

$ws.Range("C3").Value = "C00005"
$ws.Range("C4").Value = "C00006"

$ws.Range("D3").Value = "111"
$ws.Range("D4").Value = "222"

$ws.Range("E3").Value = "9001-1234"
$ws.Range("E4").Value = "9001-0004"

$ws.Range("F3").Value = "JPM2"
$ws.Range("F4").Value = "JPM2"

$ws.Range("G3").Value = "REDEMPTION"
$ws.Range("G4").Value = "REDEMPTION"

$ws.Range("H3").Value = 10
$ws.Range("H4").Value = 20

$ws.Range("I3").Value = "HKD"
$ws.Range("I4").Value = "HKD"

$ws.Range("J3").Value = "{current_timestamp}"
$ws.Range("J4").Value = "{current_timestamp}"

$ws.Range("K3").Value = "JPM2"
$ws.Range("K4").Value = "JPM2"

$ws.Range("L3").Formula = "=H3+4000"
$ws.Range("L4").Formula = "=H4+3000"

# --- 4. Row 2 gets a real formula instead of the hard-coded 1000 ---------
$ws.Range("L2").Formula = "=H2+3000"

# --- 5. Row 5 (old row 3) also switches from a static value to a formula -
$ws.Range("L5").Formula = "=H5+3000"

# --- 6. Renumber the remaining order ids (rows 5-8 used to be 2,2,2,3) --
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5

# --- 7. Final selection matches the author's last edit (cell L4) --------
$ws.Range("L4").Select()
